$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin Price (D) and Volume(1h) (E) columns to refreshed values.
# D-column values are plain text (not numbers) in the source data, so we
# force text entry (NumberFormat "@") for any numeric-looking value, then
# clear the temporary format so the cell keeps its original (unstyled) look.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.617.62'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.00%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.843.73'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.69%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.94'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.29%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4227'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.70%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3635'
$ws.Range('D8').ClearFormats()

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.32'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.86%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07256'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.06%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8892'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -5.17%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.61'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.36%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.831.62'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.57%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.562'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.87%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.328'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.81%  '

$ws.Range('E16').Value = '  -0.28%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.06%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '78.85'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.07%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008858'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.90%  '

$ws.Range('E20').Value = '  -0.18%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.44'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.79%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.608.90'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.98%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.978'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.35%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.54'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.71%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.050.62'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.59%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.965'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.94%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.98'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.59%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.52'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.08%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '120.28'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +6.14%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.243'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.94%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.850'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +6.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08898'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.79%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7785'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.67%  '

$ws.Range('E34').Value = '  -5.29%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.937'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.096'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.29%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.0000'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.19%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05395'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.35%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.099'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.42%  '

$ws.Range('E40').Value = '  -1.82%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.799'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.48%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.864'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.84%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5070'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.17%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1650'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.71%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.273'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -5.33%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.06600'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.62%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.34'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.41%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4702'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.48%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.46'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.16%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.000'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.632'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.57%  '
